$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 234.75
$ws.Range("I6").Value = 204.66667
$ws.Range("K6").Value = 614.00001
$ws.Range("M6").Value = -502.00001

$ws.Range("H29").Value = 1127.2727
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 2750
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 8250
$ws.Range("M29").Value = -319
$ws.Range("N29").Value = -8812

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H70").Value = 2180
$ws.Range("I70").Value = 1966.6666
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 5899.9998
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -5629.9998
$ws.Range("N70").Value = -8040

$ws.Range("H73").Value = 2180
$ws.Range("I73").Value = 1966.6666
$ws.Range("J73").Value = 2500
$ws.Range("K73").Value = 5899.9998
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -4963.9998
$ws.Range("N73").Value = -9372

$ws.Range("H97").Value = 8028.5
$ws.Range("J97").Value = 8028.5
$ws.Range("L97").Value = 24085.5
$ws.Range("N97").Value = -25077.5

$ws.Range("H112").Value = 7874.25
$ws.Range("J112").Value = 9999
$ws.Range("L112").Value = 29997
$ws.Range("N112").Value = -32213

$ws.Range("H134").Value = 80000
$ws.Range("J134").Value = 80000
$ws.Range("L134").Value = 80000
$ws.Range("N134").Value = -90140

$ws.Range("H137").Value = 4800.125
$ws.Range("I137").Value = 4721.6
$ws.Range("K137").Value = 14164.8
$ws.Range("M137").Value = -11614.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H44").Value = 31000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35976

$ws.Range("H45").Value = 2618.2856
$ws.Range("I45").Value = 2419.697
$ws.Range("K45").Value = 2419.697
$ws.Range("M45").Value = -2042.697

$ws.Range("H61").Value = 4938.25
$ws.Range("I61").Value = 4938.25
$ws.Range("K61").Value = 4938.25
$ws.Range("M61").Value = -4726.25

$ws.Range("H74").Value = 967.38464
$ws.Range("I74").Value = 967.38464
$ws.Range("K74").Value = 967.38464
$ws.Range("M74").Value = -93.38463999999999

$ws.Range("H77").Value = 967.38464
$ws.Range("I77").Value = 967.38464
$ws.Range("K77").Value = 4836.9232
$ws.Range("M77").Value = -468.9232000000002

$ws.Range("H97").Value = 674.73334
$ws.Range("I97").Value = 543.5
$ws.Range("K97").Value = 543.5
$ws.Range("M97").Value = -47.5

$ws.Range("H136").Value = 4938.25
$ws.Range("I136").Value = 4938.25
$ws.Range("K136").Value = 14814.75
$ws.Range("M136").Value = -12264.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5157.8887
$ws.Range("I134").Value = 5052.5
$ws.Range("K134").Value = 15157.5
$ws.Range("M134").Value = -12622.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5788.6665
$ws.Range("I86").Value = 8249.166999999999
$ws.Range("J86").Value = 3328.1667
$ws.Range("K86").Value = 8249.166999999999
$ws.Range("L86").Value = 3328.1667
$ws.Range("M86").Value = -7126.166999999999
$ws.Range("N86").Value = -5574.1667

$ws.Range("H89").Value = 5788.6665
$ws.Range("I89").Value = 8249.166999999999
$ws.Range("J89").Value = 3328.1667
$ws.Range("K89").Value = 41245.835
$ws.Range("L89").Value = 16640.8335
$ws.Range("M89").Value = -35629.835
$ws.Range("N89").Value = -27872.8335

$ws.Range("H141").Value = 200000
$ws.Range("J141").Value = 200000
$ws.Range("L141").Value = 200000
$ws.Range("N141").Value = -210360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1.3333334
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 2
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 6
$ws.Range("M17").Value = 166
$ws.Range("N17").Value = -344

$ws.Range("H34").Value = 5000
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15168

$ws.Range("H39").Value = 895.6667
$ws.Range("J39").Value = 904.8
$ws.Range("L39").Value = 2714.4
$ws.Range("N39").Value = -3302.4

$ws.Range("H55").Value = 1000
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354

$ws.Range("H122").Value = 8333.333000000001
$ws.Range("I122").Value = 8000
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 72000
$ws.Range("L122").Value = 81000
$ws.Range("M122").Value = -69550
$ws.Range("N122").Value = -85900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 20000
$ws.Range("J46").Value = 20000
$ws.Range("L46").Value = 20000
$ws.Range("N46").Value = -20312

$ws.Range("H57").Value = 22333.334
$ws.Range("J57").Value = 25000
$ws.Range("L57").Value = 25000
$ws.Range("N57").Value = -26640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H62").Value = 3249.5
$ws.Range("I62").Value = 3249.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3249.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2625.5
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3249.5
$ws.Range("I65").Value = 3249.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16247.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -13127.5
$ws.Range("N65").ClearContents()

$ws.Range("H81").Value = 699
$ws.Range("I81").Value = 699
$ws.Range("K81").Value = 1398
$ws.Range("M81").Value = -337

$ws.Range("H84").Value = 699
$ws.Range("I84").Value = 699
$ws.Range("K84").Value = 6990
$ws.Range("M84").Value = -1686

$ws.Range("H107").Value = 1247.25
$ws.Range("I107").Value = 695.5714
$ws.Range("K107").Value = 2086.7142
$ws.Range("M107").Value = -166.7142000000003

$ws.Range("H132").Value = 6705.0527
$ws.Range("I132").Value = 6212.25
$ws.Range("K132").Value = 18636.75
$ws.Range("M132").Value = -16106.75
